# Updated CircadiPy cosinor results (re-run of analysis, sawtooth_10, cosinor_7)
# Mirrors: "Make figures again to publication" -- re-ran CircaDB + CircadiPy
# simulation analyses, refreshing the computed statistics in row 2 and row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 2 -----
$ws.Range("E2").Value = 22.85000000000013
$ws.Range("G2").Value = [double]"4.417321410876696e-07"
$ws.Range("H2").Value = [double]"6.902387004346253e-06"
$ws.Range("K2").Value = 5.477220717254416
$ws.Range("L2").Value = "[3.046115563607443, 7.908325870901389]"
$ws.Range("M2").Value = [double]"1.234158585528e-05"
$ws.Range("N2").Value = [double]"1.234158585528e-05"
$ws.Range("O2").Value = -1.320789704211925
$ws.Range("P2").Value = "[-1.8365266363327715, -0.8050527720910781]"
$ws.Range("Q2").Value = [double]"7.372013777207087e-07"
$ws.Range("R2").Value = [double]"7.372013777207087e-07"
$ws.Range("S2").Value = 10.7461310998081
$ws.Range("T2").Value = "[9.359953149631707, 12.132309049984485]"
$ws.Range("W2").Value = 4.803303303303331
$ws.Range("X2").Value = 2.927727727727745
$ws.Range("Y2").Value = 6.678878878878916

# ----- Row 3 -----
$ws.Range("E3").Value = 24.32000000000036
$ws.Range("G3").Value = [double]"1.026028664252721e-06"
$ws.Range("H3").Value = [double]"6.902387004346253e-06"
$ws.Range("I3").ClearContents()
$ws.Range("K3").Value = 6.11595615451308
$ws.Range("L3").Value = "[3.472694667282486, 8.759217641743675]"
$ws.Range("M3").Value = [double]"7.798571893147255e-06"
$ws.Range("N3").Value = [double]"1.234158585528e-05"
$ws.Range("O3").Value = 1.717026615475503
$ws.Range("P3").Value = "[1.2264475824825034, 2.207605648468503]"
$ws.Range("Q3").Value = [double]"3.594946562657242e-11"
$ws.Range("R3").Value = [double]"7.189893125314484e-11"
$ws.Range("S3").Value = 10.39978191795657
$ws.Range("T3").Value = "[8.81889255744547, 11.980671278467664]"
$ws.Range("W3").Value = 17.67399399399426
$ws.Range("X3").Value = 15.77513513513537
$ws.Range("Y3").Value = 19.57285285285314
